$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.248062015503876
$ws.Range("C2").Value = 0.4780361757105943
$ws.Range("J2").Value = 0.020671834625323
$ws.Range("P2").Value = 0.186046511627907
$ws.Range("S2").Value = 0.06718346253229975
$ws.Range("B3").Value = 0.01052631578947368
$ws.Range("C3").Value = 0.01578947368421053
$ws.Range("J3").Value = 0.01578947368421053
$ws.Range("P3").Value = 0.7052631578947368
$ws.Range("S3").Value = 0.2526315789473684
$ws.Range("J4").Value = 0.05263157894736842
$ws.Range("P4").Value = 0.6842105263157895
$ws.Range("S4").Value = 0.2631578947368421
$ws.Range("B6").Value = 0.04716981132075472
$ws.Range("D6").Value = 0.02358490566037736
$ws.Range("F6").Value = 0.08962264150943396
$ws.Range("J6").Value = 0.2688679245283019
$ws.Range("O6").Value = 0.009433962264150943
$ws.Range("Q6").Value = 0.1226415094339623
$ws.Range("R6").Value = 0.04716981132075472
$ws.Range("S6").Value = 0.3915094339622642
$ws.Range("B7").Value = 0.1377777777777778
$ws.Range("D7").Value = 0.02666666666666667
$ws.Range("F7").Value = 0.07111111111111111
$ws.Range("J7").Value = 0.1333333333333333
$ws.Range("O7").Value = 0.01777777777777778
$ws.Range("Q7").Value = 0.1777777777777778
$ws.Range("R7").Value = 0.07111111111111111
$ws.Range("S7").Value = 0.3644444444444445
$ws.Range("B8").Value = 0.1033755274261603
$ws.Range("D8").Value = 0.01476793248945148
$ws.Range("F8").Value = 0.0379746835443038
$ws.Range("J8").Value = 0.1350210970464135
$ws.Range("O8").Value = 0.02109704641350211
$ws.Range("Q8").Value = 0.1518987341772152
$ws.Range("R8").Value = 0.0970464135021097
$ws.Range("S8").Value = 0.4388185654008439
$ws.Range("B9").Value = 0.07514450867052024
$ws.Range("D9").Value = 0.02312138728323699
$ws.Range("F9").Value = 0.05202312138728324
$ws.Range("J9").Value = 0.115606936416185
$ws.Range("O9").Value = 0.02890173410404624
$ws.Range("Q9").Value = 0.208092485549133
$ws.Range("R9").Value = 0.09248554913294797
$ws.Range("S9").Value = 0.4046242774566474
$ws.Range("B10").Value = 0.1410727406318883
$ws.Range("D10").Value = 0.02645113886847906
$ws.Range("E10").Value = 0.001469507714915503
$ws.Range("F10").Value = 0.05363703159441587
$ws.Range("J10").Value = 0.121234386480529
$ws.Range("O10").Value = 0.01763409257898604
$ws.Range("Q10").Value = 0.2211609110947833
$ws.Range("R10").Value = 0.08890521675238795
$ws.Range("S10").Value = 0.328434974283615
$ws.Range("G11").Value = 0.1522388059701492
$ws.Range("J11").Value = 0.09552238805970149
$ws.Range("K11").Value = 0.1850746268656716
$ws.Range("L11").Value = 0.5552238805970149
$ws.Range("S11").Value = 0.01194029850746269
$ws.Range("G12").Value = 0.7821782178217822
$ws.Range("J12").Value = 0.1534653465346535
$ws.Range("K12").Value = 0.009900990099009901
$ws.Range("L12").Value = 0.0198019801980198
$ws.Range("S12").Value = 0.03465346534653466
$ws.Range("G13").Value = 0.7692307692307693
$ws.Range("J13").Value = 0.1794871794871795
$ws.Range("S13").Value = 0.05128205128205128
$ws.Range("F15").Value = 0.01339285714285714
$ws.Range("H15").Value = 0.1383928571428572
$ws.Range("I15").Value = 0.04910714285714286
$ws.Range("J15").Value = 0.3883928571428572
$ws.Range("K15").Value = 0.06696428571428571
$ws.Range("M15").Value = 0.01339285714285714
$ws.Range("O15").Value = 0.0625
$ws.Range("S15").Value = 0.2678571428571428
$ws.Range("F16").Value = 0.02892561983471074
$ws.Range("H16").Value = 0.140495867768595
$ws.Range("I16").Value = 0.08677685950413223
$ws.Range("J16").Value = 0.3966942148760331
$ws.Range("K16").Value = 0.128099173553719
$ws.Range("M16").Value = 0.008264462809917356
$ws.Range("O16").Value = 0.06611570247933884
$ws.Range("S16").Value = 0.1446280991735537
$ws.Range("F17").Value = 0.0235042735042735
$ws.Range("H17").Value = 0.1495726495726496
$ws.Range("I17").Value = 0.07051282051282051
$ws.Range("J17").Value = 0.4423076923076923
$ws.Range("K17").Value = 0.08547008547008547
$ws.Range("M17").Value = 0.01282051282051282
$ws.Range("O17").Value = 0.07692307692307693
$ws.Range("S17").Value = 0.1388888888888889
$ws.Range("F18").Value = 0.01923076923076923
$ws.Range("H18").Value = 0.1730769230769231
$ws.Range("I18").Value = 0.08173076923076923
$ws.Range("J18").Value = 0.4278846153846154
$ws.Range("K18").Value = 0.08653846153846154
$ws.Range("M18").Value = 0.01923076923076923
$ws.Range("N18").Value = 0.004807692307692308
$ws.Range("O18").Value = 0.0576923076923077
$ws.Range("S18").Value = 0.1298076923076923
$ws.Range("F19").Value = 0.01933085501858736
$ws.Range("H19").Value = 0.2245353159851301
$ws.Range("I19").Value = 0.06617100371747212
$ws.Range("J19").Value = 0.3576208178438662
$ws.Range("K19").Value = 0.1174721189591078
$ws.Range("M19").Value = 0.0171003717472119
$ws.Range("O19").Value = 0.05427509293680297
$ws.Range("S19").Value = 0.1434944237918216
